# Auto-generated market-price update for Gilgamesh leve profit tracker.
# Mirrors the scheduled runner's refresh of currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) pulled from the Universalis market board API.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 320.8095
$ws.Range("I33").Value = 334.05
$ws.Range("K33").Value = 334.05
$ws.Range("M33").Value = -105.05
$ws.Range("H74").Value = 12888.772
$ws.Range("I74").Value = 13026.381
$ws.Range("K74").Value = 13026.381
$ws.Range("M74").Value = -12090.381
$ws.Range("H77").Value = 12888.772
$ws.Range("I77").Value = 13026.381
$ws.Range("K77").Value = 65131.905
$ws.Range("M77").Value = -60451.905
$ws.Range("H86").Value = 4801.222
$ws.Range("I86").Value = 4812.0713
$ws.Range("J86").Value = 4763.25
$ws.Range("K86").Value = 4812.0713
$ws.Range("L86").Value = 4763.25
$ws.Range("M86").Value = -3689.0713
$ws.Range("N86").Value = -7009.25
$ws.Range("H89").Value = 4801.222
$ws.Range("I89").Value = 4812.0713
$ws.Range("J89").Value = 4763.25
$ws.Range("K89").Value = 24060.3565
$ws.Range("L89").Value = 23816.25
$ws.Range("M89").Value = -18444.3565
$ws.Range("N89").Value = -35048.25
$ws.Range("H98").Value = 2132.077
$ws.Range("I98").Value = 2132.077
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2132.077
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -634.0770000000002
$ws.Range("H103").Value = 3649.7
$ws.Range("H107").Value = 341.84616
$ws.Range("I107").Value = 355.16666
$ws.Range("K107").Value = 355.16666
$ws.Range("M107").Value = 1564.83334
$ws.Range("H112").Value = 2021.4286
$ws.Range("J112").Value = 2210.25
$ws.Range("L112").Value = 6630.75
$ws.Range("N112").Value = -8846.75
$ws.Range("H115").Value = 817
$ws.Range("I115").Value = 423.66666
$ws.Range("K115").Value = 1270.99998
$ws.Range("M115").Value = 296.0000199999999
$ws.Range("H121").Value = 1814.2142
$ws.Range("J121").Value = 1842
$ws.Range("L121").Value = 5526
$ws.Range("N121").Value = -9020
$ws.Range("H122").Value = 2132.077
$ws.Range("I122").Value = 2132.077
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6396.231000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3946.231000000001
$ws.Range("H135").Value = 1515.9445
$ws.Range("I135").Value = 1385.8667
$ws.Range("K135").Value = 12472.8003
$ws.Range("M135").Value = -9937.800300000001
$ws.Range("H137").Value = 2006939
$ws.Range("I137").Value = 2501123.5
$ws.Range("K137").Value = 7503370.5
$ws.Range("M137").Value = -7500820.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2067432.4
$ws.Range("I32").Value = 1031665.7
$ws.Range("K32").Value = 1031665.7
$ws.Range("M32").Value = -1031378.7
$ws.Range("H61").Value = 4301.3
$ws.Range("I61").Value = 8006.5
$ws.Range("J61").Value = 3375
$ws.Range("K61").Value = 8006.5
$ws.Range("L61").Value = 3375
$ws.Range("M61").Value = -7794.5
$ws.Range("N61").Value = -3799
$ws.Range("H74").Value = 15713886
$ws.Range("I74").Value = 137045.7
$ws.Range("K74").Value = 137045.7
$ws.Range("M74").Value = -136171.7
$ws.Range("H77").Value = 15713886
$ws.Range("I77").Value = 137045.7
$ws.Range("K77").Value = 685228.5
$ws.Range("M77").Value = -680860.5
$ws.Range("H97").Value = 2079.7778
$ws.Range("I97").Value = 1796
$ws.Range("J97").Value = 4350
$ws.Range("K97").Value = 1796
$ws.Range("L97").Value = 4350
$ws.Range("M97").Value = -1300
$ws.Range("N97").Value = -5342
$ws.Range("H132").Value = 1674149.4
$ws.Range("I132").Value = 2565713.5
$ws.Range("K132").Value = 7697140.5
$ws.Range("M132").Value = -7694610.5
$ws.Range("H136").Value = 4301.3
$ws.Range("I136").Value = 8006.5
$ws.Range("J136").Value = 3375
$ws.Range("K136").Value = 24019.5
$ws.Range("L136").Value = 10125
$ws.Range("M136").Value = -21469.5
$ws.Range("N136").Value = -15225

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5001.5557
$ws.Range("I99").Value = 4430.5713
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 4430.5713
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = -2932.5713
$ws.Range("N99").Value = -9996
$ws.Range("H105").Value = 13686792
$ws.Range("I105").Value = 1001889.1
$ws.Range("K105").Value = 1001889.1
$ws.Range("M105").Value = -1000142.1
$ws.Range("H134").Value = 2587.1667
$ws.Range("I134").Value = 2024.6666
$ws.Range("K134").Value = 6073.9998
$ws.Range("M134").Value = -3538.9998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1926075.1
$ws.Range("I31").Value = 3500
$ws.Range("J31").Value = 2019103
$ws.Range("K31").Value = 3500
$ws.Range("L31").Value = 2019103
$ws.Range("M31").Value = -3205
$ws.Range("N31").Value = -2019693
$ws.Range("H34").Value = 1926075.1
$ws.Range("I34").Value = 3500
$ws.Range("J34").Value = 2019103
$ws.Range("K34").Value = 3500
$ws.Range("L34").Value = 2019103
$ws.Range("M34").Value = -3298
$ws.Range("N34").Value = -2019507
$ws.Range("H58").Value = 5547.478
$ws.Range("I58").Value = 3608.7
$ws.Range("J58").Value = 7038.846
$ws.Range("K58").Value = 3608.7
$ws.Range("L58").Value = 7038.846
$ws.Range("M58").Value = -3405.7
$ws.Range("N58").Value = -7444.846
$ws.Range("H132").Value = 2993.55
$ws.Range("I132").Value = 2848.9644
$ws.Range("J132").Value = 3330.9167
$ws.Range("K132").Value = 8546.893199999999
$ws.Range("L132").Value = 9992.750100000001
$ws.Range("M132").Value = -6016.893199999999
$ws.Range("N132").Value = -15052.7501
$ws.Range("H134").Value = 3496.9333
$ws.Range("I134").Value = 3424.238
$ws.Range("J134").Value = 3666.5557
$ws.Range("K134").Value = 10272.714
$ws.Range("L134").Value = 10999.6671
$ws.Range("M134").Value = -7737.714
$ws.Range("N134").Value = -16069.6671
$ws.Range("H136").Value = 5547.478
$ws.Range("I136").Value = 3608.7
$ws.Range("J136").Value = 7038.846
$ws.Range("K136").Value = 10826.1
$ws.Range("L136").Value = 21116.538
$ws.Range("M136").Value = -8276.099999999999
$ws.Range("N136").Value = -26216.538

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 464.92856
$ws.Range("I5").Value = 427.77777
$ws.Range("K5").Value = 1283.33331
$ws.Range("M5").Value = -1171.33331
$ws.Range("H33").Value = 20.5
$ws.Range("J33").Value = 20.5
$ws.Range("L33").Value = 123
$ws.Range("N33").Value = -689
$ws.Range("H56").Value = 7764.5625
$ws.Range("I56").Value = 7764.5625
$ws.Range("K56").Value = 7764.5625
$ws.Range("M56").Value = -7234.5625
$ws.Range("H64").Value = 2004.5454
$ws.Range("J64").Value = 2155
$ws.Range("L64").Value = 6465
$ws.Range("N64").Value = -7005
$ws.Range("H67").Value = 2004.5454
$ws.Range("J67").Value = 2155
$ws.Range("L67").Value = 6465
$ws.Range("N67").Value = -8337
$ws.Range("H114").Value = 1904.4375
$ws.Range("I114").Value = 785.5714
$ws.Range("J114").Value = 2774.6667
$ws.Range("K114").Value = 2356.7142
$ws.Range("L114").Value = 8324.000100000001
$ws.Range("M114").Value = 897.2857999999997
$ws.Range("N114").Value = -14832.0001
$ws.Range("H117").Value = 2756.3333
$ws.Range("I117").Value = 2500
$ws.Range("J117").Value = 2884.5
$ws.Range("K117").Value = 7500
$ws.Range("L117").Value = 8653.5
$ws.Range("M117").Value = -4058
$ws.Range("N117").Value = -15537.5
$ws.Range("H120").Value = 31499.75
$ws.Range("J120").Value = 33000
$ws.Range("L120").Value = 99000
$ws.Range("N120").Value = -108676
$ws.Range("H132").Value = 5213.8086
$ws.Range("I132").Value = 5376.727
$ws.Range("J132").Value = 5164.028
$ws.Range("K132").Value = 48390.543
$ws.Range("L132").Value = 46476.252
$ws.Range("M132").Value = -45860.543
$ws.Range("N132").Value = -51536.252
$ws.Range("H135").Value = 464.92856
$ws.Range("I135").Value = 427.77777
$ws.Range("K135").Value = 3849.99993
$ws.Range("M135").Value = -1314.99993

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 991.6667
$ws.Range("I46").Value = 890
$ws.Range("K46").Value = 890
$ws.Range("M46").Value = -702
$ws.Range("H93").Value = 2289.6
$ws.Range("I93").Value = 2430.625
$ws.Range("J93").Value = 2128.4285
$ws.Range("K93").Value = 2430.625
$ws.Range("L93").Value = 2128.4285
$ws.Range("M93").Value = -1182.625
$ws.Range("N93").Value = -4624.4285
$ws.Range("H128").Value = 44999
$ws.Range("J128").Value = 44999
$ws.Range("L128").Value = 44999
$ws.Range("N128").Value = -54959
$ws.Range("H132").Value = 6799.2144
$ws.Range("I132").Value = 16330.667
$ws.Range("J132").Value = 4199.727
$ws.Range("K132").Value = 48992.001
$ws.Range("L132").Value = 12599.181
$ws.Range("M132").Value = -46462.001
$ws.Range("N132").Value = -17659.181
$ws.Range("H136").Value = 7353.3335
$ws.Range("I136").Value = 5619.3335
$ws.Range("K136").Value = 16858.0005
$ws.Range("M136").Value = -14308.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 57999.5
$ws.Range("J46").Value = 57999.5
$ws.Range("L46").Value = 57999.5
$ws.Range("N46").Value = -58461.5
$ws.Range("H132").Value = 2365.432
$ws.Range("I132").Value = 2299.5945
$ws.Range("J132").Value = 2713.4285
$ws.Range("K132").Value = 6898.7835
$ws.Range("L132").Value = 8140.2855
$ws.Range("M132").Value = -4368.7835
$ws.Range("N132").Value = -13200.2855
$ws.Range("H134").Value = 57999.5
$ws.Range("J134").Value = 57999.5
$ws.Range("L134").Value = 173998.5
$ws.Range("N134").Value = -179068.5
$ws.Range("H136").Value = 8555993
$ws.Range("I136").Value = 12830062
$ws.Range("J136").Value = 7854
$ws.Range("K136").Value = 38490186
$ws.Range("L136").Value = 23562
$ws.Range("M136").Value = -38487636
$ws.Range("N136").Value = -28662

